# "Committed Corporate Customer excel file"
# Adds the Corporate-Customer TDR fields between PRINCIPAL and AUTO.ROLLOVER,
# moving the AUTO.ROLLOVER column out to the right of the new block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the additional Corporate Customer columns, inserted
# between the existing PRINCIPAL (D) and AUTO.ROLLOVER columns.
$ws.Range("E1").Value = "INTEND.DATE"
$ws.Range("F1").Value = "CUST.REMARKS:1"
$ws.Range("G1").Value = "TAX.INTEREST.TYPE:1"
$ws.Range("H1").Value = "DRAWDOWN.ACCOUNT"
$ws.Range("I1").Value = "PRIN.LIQ.ACCT"
$ws.Range("J1").Value = "INT.LIQ.ACCT"
$ws.Range("K1").Value = "CHRG.LIQ.ACCT"

# AUTO.ROLLOVER now lives in column L, after the new block.
$ws.Range("L1").Value = "AUTO.ROLLOVER"
$ws.Range("M1").Value = "FINAL.MATURITY"
$ws.Range("N1").Value = "EXP.DATE"

# Move the AUTO.ROLLOVER data value from its old spot (E2) to the new one (L2).
$ws.Range("E2").Value = $null
$ws.Range("L2").Value = 2

# Columns E:K take the narrower standard width used throughout the sheet,
# while the column that now holds AUTO.ROLLOVER (L) keeps the wider width
# that used to belong to the old E column.
$ws.Range("E1:K1").ColumnWidth = 9.5
$ws.Range("L1:L1").ColumnWidth = 15

# Restore the previously-minimized workbook window.
$excel.ActiveWindow.WindowState = -4143

# Selection moves to D15, with N1 also part of the selected set.
$ws.Activate() | Out-Null
$ws.Range("D15,N1").Select() | Out-Null
